# Update iServ stats for 2025-11 (row 24)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B24").Value = 6327
$ws.Range("D24").Value = 5924488
$ws.Range("E24").Value = 936.3818555397503
$ws.Range("F24").Value = 7.858847596317764
$ws.Range("H24").Value = 25.50360422423383
